$wb = $excel.ActiveWorkbook

# Add the new worksheet "توزین" right after Sheet1 (matches final tab order
# Sheet1, توزین and the saved activeTab pointing at the new sheet).
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "توزین"

# Title row
$ws2.Range("A1").Value = "بر آورد وزن اشیا"

# Header row
$ws2.Range("A2").Value = "ردیف"
$ws2.Range("B2").Value = "نام"
$ws2.Range("C2").Value = "وزن"
$ws2.Range("D2").Value = "وزن کل"

# First data row
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "چمدان بزرگ"
$ws2.Range("C3").Value = 5
$ws2.Range("D3").Formula = "=C3"

# Running-total column filled down through row 23 (shared formula)
$ws2.Range("D4:D23").Formula = "=C4+C3"

# Auto-fit column A like the original author did
$ws2.Columns.Item(1).AutoFit()

# Leave the selection where the author left it
$ws2.Range("D4").Select()
